$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 12500519
$ws.Range("I103").Value = 344.5
$ws.Range("J103").Value = 20000624
$ws.Range("K103").Value = 1033.5
$ws.Range("L103").Value = 60001872
$ws.Range("M103").Value = -447.5
$ws.Range("N103").Value = -60003044

$ws.Range("H132").Value = 2132.4426
$ws.Range("I132").Value = 2049.2827
$ws.Range("J132").Value = 2387.4666
$ws.Range("K132").Value = 6147.848100000001
$ws.Range("L132").Value = 7162.399800000001
$ws.Range("M132").Value = -3617.848100000001
$ws.Range("N132").Value = -12222.3998

$ws.Range("H134").Value = 41597.7
$ws.Range("J134").Value = 41597.7
$ws.Range("L134").Value = 41597.7
$ws.Range("N134").Value = -51737.7

$ws.Range("H137").Value = 3001229.5
$ws.Range("I137").Value = 1220658
$ws.Range("J137").Value = 11112723
$ws.Range("K137").Value = 3661974
$ws.Range("L137").Value = 33338169
$ws.Range("M137").Value = -3659424
$ws.Range("N137").Value = -33343269

$ws.Range("H138").Value = 278960.34
$ws.Range("I138").Value = 979.8333
$ws.Range("J138").Value = 2224824
$ws.Range("K138").Value = 2939.4999
$ws.Range("L138").Value = 6674472
$ws.Range("M138").Value = 2200.5001
$ws.Range("N138").Value = -6684752

$ws.Range("H141").Value = 4275725.5
$ws.Range("I141").Value = 2156.7727
$ws.Range("J141").Value = 9806227
$ws.Range("K141").Value = 6470.3181
$ws.Range("L141").Value = 29418681
$ws.Range("M141").Value = -1290.3181
$ws.Range("N141").Value = -29429041

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4890.54
$ws.Range("I32").Value = 3817.484
$ws.Range("J32").Value = 19146.857
$ws.Range("K32").Value = 3817.484
$ws.Range("L32").Value = 19146.857
$ws.Range("M32").Value = -3530.484
$ws.Range("N32").Value = -19720.857

$ws.Range("H61").Value = 1306.6522
$ws.Range("I61").Value = 1306.0952
$ws.Range("J61").Value = 1312.5
$ws.Range("K61").Value = 1306.0952
$ws.Range("L61").Value = 1312.5
$ws.Range("M61").Value = -1094.0952
$ws.Range("N61").Value = -1736.5

$ws.Range("H74").Value = 857.95654
$ws.Range("I74").Value = 846.13336
$ws.Range("J74").Value = 1390
$ws.Range("K74").Value = 846.13336
$ws.Range("L74").Value = 1390
$ws.Range("M74").Value = 27.86663999999996
$ws.Range("N74").Value = -3138

$ws.Range("H77").Value = 857.95654
$ws.Range("I77").Value = 846.13336
$ws.Range("J77").Value = 1390
$ws.Range("K77").Value = 4230.6668
$ws.Range("L77").Value = 6950
$ws.Range("M77").Value = 137.3332
$ws.Range("N77").Value = -15686

$ws.Range("H97").Value = 643.3182
$ws.Range("I97").Value = 622.86487
$ws.Range("J97").Value = 751.4286
$ws.Range("K97").Value = 622.86487
$ws.Range("L97").Value = 751.4286
$ws.Range("M97").Value = -126.86487
$ws.Range("N97").Value = -1743.4286

$ws.Range("H124").Value = 21732.334
$ws.Range("J124").Value = 21732.334
$ws.Range("L124").Value = 21732.334
$ws.Range("N124").Value = -31552.334

$ws.Range("H136").Value = 1306.6522
$ws.Range("I136").Value = 1306.0952
$ws.Range("J136").Value = 1312.5
$ws.Range("K136").Value = 3918.2856
$ws.Range("L136").Value = 3937.5
$ws.Range("M136").Value = -1368.2856
$ws.Range("N136").Value = -9037.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 398.90475
$ws.Range("I80").Value = 189.16667
$ws.Range("J80").Value = 482.8
$ws.Range("K80").Value = 189.16667
$ws.Range("L80").Value = 482.8
$ws.Range("M80").Value = 808.8333299999999
$ws.Range("N80").Value = -2478.8

$ws.Range("H83").Value = 398.90475
$ws.Range("I83").Value = 189.16667
$ws.Range("J83").Value = 482.8
$ws.Range("K83").Value = 945.8333500000001
$ws.Range("L83").Value = 2414
$ws.Range("M83").Value = 4046.16665
$ws.Range("N83").Value = -12398

$ws.Range("H134").Value = 98716.39
$ws.Range("I134").Value = 138030.9
$ws.Range("J134").Value = 2614.2222
$ws.Range("K134").Value = 414092.7
$ws.Range("L134").Value = 7842.6666
$ws.Range("M134").Value = -411557.7
$ws.Range("N134").Value = -12912.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4493
$ws.Range("I16").Value = 4493
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 4493
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -4206
$ws.Range("N16").ClearContents()

$ws.Range("H31").Value = 1533.8889
$ws.Range("I31").Value = 1125.5555
$ws.Range("J31").Value = 3575.5557
$ws.Range("K31").Value = 1125.5555
$ws.Range("L31").Value = 3575.5557
$ws.Range("M31").Value = -830.5554999999999
$ws.Range("N31").Value = -4165.5557

$ws.Range("H34").Value = 1533.8889
$ws.Range("I34").Value = 1125.5555
$ws.Range("J34").Value = 3575.5557
$ws.Range("K34").Value = 1125.5555
$ws.Range("L34").Value = 3575.5557
$ws.Range("M34").Value = -923.5554999999999
$ws.Range("N34").Value = -3979.5557

$ws.Range("H58").Value = 1209.919
$ws.Range("I58").Value = 1244
$ws.Range("J58").Value = 991.8
$ws.Range("K58").Value = 1244
$ws.Range("L58").Value = 991.8
$ws.Range("M58").Value = -1041
$ws.Range("N58").Value = -1397.8

$ws.Range("H113").Value = 4493
$ws.Range("I113").Value = 4493
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4493
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -2323
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 2891.7778
$ws.Range("I132").Value = 2350.077
$ws.Range("J132").Value = 4300.2
$ws.Range("K132").Value = 7050.231000000001
$ws.Range("L132").Value = 12900.6
$ws.Range("M132").Value = -4520.231000000001
$ws.Range("N132").Value = -17960.6

$ws.Range("H134").Value = 4450.45
$ws.Range("I134").Value = 6150.5454
$ws.Range("J134").Value = 2372.5557
$ws.Range("K134").Value = 18451.6362
$ws.Range("L134").Value = 7117.6671
$ws.Range("M134").Value = -15916.6362
$ws.Range("N134").Value = -12187.6671

$ws.Range("H136").Value = 1209.919
$ws.Range("I136").Value = 1244
$ws.Range("J136").Value = 991.8
$ws.Range("K136").Value = 3732
$ws.Range("L136").Value = 2975.4
$ws.Range("M136").Value = -1182
$ws.Range("N136").Value = -8075.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 931.9299999999999
$ws.Range("I131").Value = 532
$ws.Range("J131").Value = 952.97894
$ws.Range("K131").Value = 1596
$ws.Range("L131").Value = 2858.93682
$ws.Range("M131").Value = 3444
$ws.Range("N131").Value = -12938.93682

$ws.Range("H137").Value = 2387.8572
$ws.Range("I137").Value = 2225.5557
$ws.Range("J137").Value = 2680
$ws.Range("K137").Value = 6676.6671
$ws.Range("L137").Value = 8040
$ws.Range("M137").Value = -1576.6671
$ws.Range("N137").Value = -18240

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 3500
$ws.Range("J38").Value = 3500
$ws.Range("L38").Value = 3500
$ws.Range("N38").Value = -4426

$ws.Range("H102").Value = 1873.0294
$ws.Range("I102").Value = 1335.6842
$ws.Range("J102").Value = 2553.6667
$ws.Range("K102").Value = 1335.6842
$ws.Range("L102").Value = 2553.6667
$ws.Range("M102").Value = 286.3158000000001
$ws.Range("N102").Value = -5797.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H118").Value = 31809.2
$ws.Range("J118").Value = 31809.2
$ws.Range("L118").Value = 31809.2
$ws.Range("N118").Value = -35123.2

$ws.Range("H132").Value = 2986.3125
$ws.Range("I132").Value = 2744.96
$ws.Range("J132").Value = 3848.2856
$ws.Range("K132").Value = 8234.880000000001
$ws.Range("L132").Value = 11544.8568
$ws.Range("M132").Value = -5704.880000000001
$ws.Range("N132").Value = -16604.8568

$ws.Range("H136").Value = 1891.8857
$ws.Range("I136").Value = 1747.2
$ws.Range("J136").Value = 2760
$ws.Range("K136").Value = 5241.6
$ws.Range("L136").Value = 8280
$ws.Range("M136").Value = -2691.6
$ws.Range("N136").Value = -13380

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1809.8209
$ws.Range("I132").Value = 2052.413
$ws.Range("J132").Value = 1278.4286
$ws.Range("K132").Value = 6157.239
$ws.Range("L132").Value = 3835.2858
$ws.Range("M132").Value = -3627.239
$ws.Range("N132").Value = -8895.2858
